$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 6 values that repeat in each block of rows (offsets 0..5 within the block)
$newValues = @(3164602900, 3176652286, 3173506184, 3172486789, 3173809096, 3123144985)

# Blocks start at row 5, then every 10 rows thereafter, up to row 315 (315..320), mirroring the
# repeating pattern already present in the sheet (rows 1-4 are a fixed header block that is
# untouched, followed by 6 rows that get updated, repeated all the way to row 320).
for ($blockStart = 5; $blockStart -le 315; $blockStart += 10) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        $row = $blockStart + $i
        $ws.Cells.Item($row, 1).Value = $newValues[$i]
    }
}

# Update the view state: scrolled position and active selection
$excel.ActiveWindow.ScrollRow = 298
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A320").Select()
